$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- Resize table / columns (values are in dxa/twips; Word COM uses points = dxa/20) ---
$rowCount = $tbl.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $tbl.Cell($r, 1).Width = 2038 / 20.0
    $tbl.Cell($r, 3).Width = 2180 / 20.0
    $tbl.Cell($r, 4).Width = 2138 / 20.0
    $tbl.Cell($r, 5).Width = 2550 / 20.0
}
$tbl.PreferredWidth = 11245 / 20.0

# --- Move "Design" from the Week 6 row to the Week 7 row ---
# Week 6 = row 7, Week 7 = row 8 (row 1 is the header row)
$week6Topic = $tbl.Cell(7, 2)
$week6Rng = $week6Topic.Range
$week6Rng.Find.Execute("Function and Class Templates, Design", $true, $false, $false, $false, $false, $true, 0, $false, "Function and Class Templates", 1) | Out-Null

$week7Topic = $tbl.Cell(8, 2)
$week7Rng = $week7Topic.Range
$week7Rng.Find.Execute("Review,", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$week7Rng.Collapse(1)
$week7Rng.InsertBefore("Design, ")
